# Applies the cryptocurrency price/volume refresh described by the commit
# "Updated cryptos list on Fri Sep  6 21:50:07 UTC 2024 with GitHub Actions".
#
# All values in columns D (Price) and E (Volume(1h)) are stored as plain text
# in the source sheet (not numbers), including ones that look numeric, and
# the "Volume(1h)" strings intentionally keep their leading/trailing spaces.
# A leading single-quote is used below for Price values that look like a
# number so Excel keeps storing them as text (quoted/literal) instead of
# silently converting them into a numeric cell value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = '53.364.97'
$ws.Range("E2").Value = '  -4.71%  '

# Row 3 (Ethereum)
$ws.Range("D3").Value = '2.186.71'
$ws.Range("E3").Value = '  -7.38%  '

# Row 4 (TetherUSD)
$ws.Range("E4").Value = '  -0.12%  '

# Row 5 (BNB)
$ws.Range("D5").Value = '''481.64'
$ws.Range("E5").Value = '  -3.58%  '

# Row 6 (Solana)
$ws.Range("D6").Value = '''123.89'
$ws.Range("E6").Value = '  -4.06%  '

# Row 7 (USDC)
$ws.Range("E7").Value = '  -0.26%  '

# Row 8 (XRP)
$ws.Range("E8").Value = '  -4.82%  '

# Row 9 (LidoStakedEther)
$ws.Range("D9").Value = '2.199.81'
$ws.Range("E9").Value = '  -6.98%  '

# Row 10 (Dogecoin)
$ws.Range("D10").Value = '''0.0911'
$ws.Range("E10").Value = '  -6.98%  '

# Row 11 (TRON)
$ws.Range("E11").Value = '  -1.94%  '

# Row 12 (Cardano)
$ws.Range("E12").Value = '  -3.61%  '

# Row 13 (Toncoin)
$ws.Range("D13").Value = '''4.55'
$ws.Range("E13").Value = '  -4.91%  '

# Row 14 (WrappedliquidstakedEther2.0)
$ws.Range("D14").Value = '2.571.92'
$ws.Range("E14").Value = '  -7.56%  '

# Row 15 (Avalanche)
$ws.Range("D15").Value = '''20.90'
$ws.Range("E15").Value = '  -2.27%  '

# Row 16 (WrappedBTC)
$ws.Range("D16").Value = '53.285.40'
$ws.Range("E16").Value = '  -4.82%  '

# Row 17 (ShibaInu)
$ws.Range("E17").Value = '  -3.94%  '

# Row 18 (WrappedEther)
$ws.Range("D18").Value = '2.203.78'
$ws.Range("E18").Value = '  -8.08%  '

# Row 19 (Chainlink)
$ws.Range("E19").Value = '  -4.89%  '

# Row 20 (Polkadot)
$ws.Range("E20").Value = '  -2.89%  '

# Row 21 (BitcoinCash)
$ws.Range("D21").Value = '''292.63'
$ws.Range("E21").Value = '  -4.63%  '

# Row 22 (Uniswap)
$ws.Range("E22").Value = '  -3.61%  '

# Row 23 (Dai)
$ws.Range("D23").Value = '''0.998'
$ws.Range("E23").Value = '  -0.22%  '

# Row 24 (Litecoin)
$ws.Range("D24").Value = '''62.41'
$ws.Range("E24").Value = '  -4.65%  '

# Row 25 (Binance-PegBSC-USD)
$ws.Range("D25").Value = '''0.994'
$ws.Range("E25").Value = '  -0.39%  '

# Row 26 (Polygon)
$ws.Range("D26").Value = '''0.363'
$ws.Range("E26").Value = '  -1.48%  '

# Row 27 (WrappedeETH)
$ws.Range("D27").Value = '2.286.87'
$ws.Range("E27").Value = '  -7.55%  '

# Row 28 (Kaspa)
$ws.Range("E28").Value = '  -2.14%  '

# Row 29 (InternetComputer(DFINITY))
$ws.Range("E29").Value = '  -3.95%  '

# Row 30 (Monero)
$ws.Range("D30").Value = '''164.81'
$ws.Range("E30").Value = '  -3.79%  '

# Row 31 (USDe)
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '''1.57'
$ws.Range("E31").Value = '  -4.28%  '

# Row 32 (PancakeSwap)
$ws.Range("B32").Value = 'USDe'
$ws.Range("C32").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D32").Value = '''0.998'
$ws.Range("E32").Value = '  -0.18%  '

# Row 33 (FirstDigitalUSD)
$ws.Range("D33").Value = '''0.993'
$ws.Range("E33").Value = '  -0.50%  '

# Row 34 (PEPE)
$ws.Range("E34").Value = '  -7.65%  '

# Row 35 (Aptos)
$ws.Range("D35").Value = '''5.65'
$ws.Range("E35").Value = '  -1.63%  '

# Row 36 (Fetch.AI)
$ws.Range("E36").Value = '  -3.29%  '

# Row 37 (EthereumClassic)
$ws.Range("D37").Value = '''17.21'
$ws.Range("E37").Value = '  -2.14%  '

# Row 38 (ImmutableX)
$ws.Range("E38").Value = '  -2.76%  '

# Row 39 (SuiNetwork)
$ws.Range("D39").Value = '''0.814'
$ws.Range("E39").Value = '  +3.22%  '

# Row 40 (OKB)
$ws.Range("D40").Value = '''35.65'
$ws.Range("E40").Value = '  -1.23%  '

# Row 41 (NEARProtocol)
$ws.Range("E41").Value = '  -5.71%  '

# Row 42 (PolygonEcosystemToken)
$ws.Range("E42").Value = '  -1.42%  '

# Row 43 (Stacks)
$ws.Range("E43").Value = '  -2.17%  '

# Row 44 (Filecoin)
$ws.Range("E44").Value = '  -3.95%  '

# Row 45 (RenderToken)
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '''122.74'
$ws.Range("E45").Value = '  -4.72%  '

# Row 46 (Aave)
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").Value = '''4.69'
$ws.Range("E46").Value = '  -0.38%  '

# Row 47 (Stellar)
$ws.Range("E47").Value = '  -3.09%  '

# Row 48 (Mantle)
$ws.Range("E48").Value = '  -6.18%  '

# Row 49 (Hedera)
$ws.Range("D49").Value = '''0.0466'
$ws.Range("E49").Value = '  -3.07%  '

# Row 50 (Bittensor)
$ws.Range("D50").Value = '''227.25'
$ws.Range("E50").Value = '  -4.96%  '

# Row 51 (VeChain)
$ws.Range("E51").Value = '  -4.12%  '
